$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F6").Value = 2810
$ws.Range("F10").Value = 374
$ws.Range("F12").Value = 313
$ws.Range("F14").Value = 5922
$ws.Range("F15").Value = 629
$ws.Range("F16").Value = 1044
$ws.Range("F17").Value = 6
$ws.Range("F18").Value = 102
$ws.Range("F21").Value = 526
$ws.Range("F23").Value = 20
$ws.Range("F24").Value = 59
$ws.Range("F25").Value = 1295
$ws.Range("F28").Value = 32
$ws.Range("F29").Value = 2051
$ws.Range("F30").Value = 170
$ws.Range("F31").Value = 348
$ws.Range("F33").Value = 3268

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F7").Value = 344
$ws.Range("F24").Value = 4047
$ws.Range("F28").Value = 129
$ws.Range("F32").Value = 192
$ws.Range("F33").Value = 17

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 1803
$ws.Range("F3").Value = 86
$ws.Range("F6").Value = 1136
$ws.Range("F8").Value = 1480

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 1803
$ws.Range("F3").Value = 86
$ws.Range("F6").Value = 1136
$ws.Range("F7").Value = 1480
$ws.Range("F13").Value = 2810
$ws.Range("F17").Value = 374
$ws.Range("F20").Value = 313
$ws.Range("F22").Value = 5922
$ws.Range("F23").Value = 629
$ws.Range("F24").Value = 1044
$ws.Range("F25").Value = 102
$ws.Range("F28").Value = 526
$ws.Range("F36").Value = 1295
$ws.Range("F39").Value = 129
$ws.Range("F41").Value = 32
$ws.Range("F44").Value = 2051
$ws.Range("F45").Value = 192
$ws.Range("F47").Value = 170
$ws.Range("F48").Value = 348
$ws.Range("F50").Value = 3268
